$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45
$ws.Range("I45").Value = 0.2284633975843539
$ws.Range("J45").Value = 0.01855976243503714

# Row 46
$ws.Range("H46").Value = 0.08028600715190851
$ws.Range("I46").Value = -0.1296176279974082

# Row 47
$ws.Range("G47").Value = -0.07715998185224648
$ws.Range("H47").Value = -0.2870636170015632

# Row 48
$ws.Range("F48").Value = 0.4234994746738243
$ws.Range("G48").Value = 0.2135958395245076

# Row 49
$ws.Range("E49").Value = 0.1431415941383551
$ws.Range("F49").Value = -0.06676204101096155

# Row 50
$ws.Range("D50").Value = 0.3151164519833668
$ws.Range("E50").Value = 0.1052128168340501

# Row 51
$ws.Range("C51").Value = 0.009253912237035311
$ws.Range("D51").Value = -0.2006497229122814

# Row 52
$ws.Range("B52").Value = 0.6215838649243215
$ws.Range("C52").Value = 0.4116802297750048

# Row 53
$ws.Range("B53").Value = -0.2766911554241067
